$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.954.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.55%  "
$ws.Range("D3").Value = "'3.301.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'409.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.65%  "
$ws.Range("D6").Value = "'112.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.09%  "
$ws.Range("D7").Value = "'3.297.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  -4.67%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.618"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").Value = "'0.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.86%  "
$ws.Range("D12").Value = "'38.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "'3.838.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "'8.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "'18.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "'3.330.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "'60.872.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.87%  "
$ws.Range("D19").Value = "'0.982"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").Value = "'10.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").Value = "'3.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("D23").Value = "'12.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("D24").Value = "'294.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "'72.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").Value = "'28.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.85%  "
$ws.Range("D28").Value = "'4.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").Value = "'0.172"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").Value = "'7.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "'7.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "'11.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "'2.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.96%  "
$ws.Range("D36").Value = "'39.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'0.0475"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").Value = "'52.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'3.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("D41").Value = "'3.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("D42").Value = "'134.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "'0.119"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "'16.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.15%  "
$ws.Range("E47").Value = "  -4.64%  "
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("D49").Value = "'20.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.82%  "
$ws.Range("D50").Value = "'2.104.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'3.641.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "
